$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 0.9451964172203481
$ws.Range("E2").Value = 0.9451964172203481

# Row 3
$ws.Range("D3").Value = 0.198670621546177
$ws.Range("E3").Value = 0.198670621546177

# Row 4
$ws.Range("D4").Value = 0.09849088316178693
$ws.Range("E4").Value = 0.09849088316178693

# Row 5
$ws.Range("D5").Value = 0.01640381964650117
$ws.Range("E5").Value = 0.01640381964650117

# Row 6
$ws.Range("D6").Value = 0.9977345164798331
$ws.Range("E6").Value = 0.9977345164798331

# Row 7
$ws.Range("C7").Value = $false
$ws.Range("D7").Value = [double]"8.038708539708431E-09"
$ws.Range("E7").Value = 0.9999999919612914

# Row 8
$ws.Range("D8").Value = 0.8855441587576345
$ws.Range("E8").Value = 0.1144558412423655

# Row 9
$ws.Range("D9").Value = 0.7476353367410732
$ws.Range("E9").Value = 0.2523646632589268

# Row 10
$ws.Range("D10").Value = 0.9999999999999076
$ws.Range("E10").Value = [double]"9.237055564881302E-14"

# Row 11
$ws.Range("D11").Value = 0.9999999999988878
$ws.Range("E11").Value = [double]"1.112221426069482E-12"
$ws.Range("F11").Value = 2.838706493377686
$ws.Range("G11").Value = 0.7

# Row 12
$ws.Range("D12").Value = 0.9305774426306326
$ws.Range("E12").Value = 0.9305774426306326

# Row 13
$ws.Range("D13").Value = 0.5946164193341191
$ws.Range("E13").Value = 0.5946164193341191

# Row 14
$ws.Range("D14").Value = 0.9999789488642293
$ws.Range("E14").Value = 0.9999789488642293

# Row 15
$ws.Range("D15").Value = 0.0003701840933868708
$ws.Range("E15").Value = 0.0003701840933868708

# Row 16
$ws.Range("D16").Value = 0.9998409436939637
$ws.Range("E16").Value = 0.9998409436939637

# Row 17
$ws.Range("C17").Value = $false
$ws.Range("D17").Value = [double]"2.514573359719146E-12"
$ws.Range("E17").Value = 0.9999999999974855

# Row 18
$ws.Range("D18").Value = 0.9990400249706943
$ws.Range("E18").Value = 0.0009599750293056708

# Row 19
$ws.Range("D19").Value = 0.9160398977756273
$ws.Range("E19").Value = 0.08396010222437267

# Row 20
$ws.Range("D20").Value = 1
$ws.Range("E20").Value = 0

# Row 21
$ws.Range("D21").Value = 1
$ws.Range("E21").Value = 0
$ws.Range("F21").Value = 4.988322257995605
$ws.Range("G21").Value = 0.5
